$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H132").Value = 23172.818
$ws.Range("I132").Value = 3713.513
$ws.Range("K132").Value = 11140.539
$ws.Range("M132").Value = -8610.539000000001
$ws.Range("H135").Value = 15625727
$ws.Range("I135").Value = 752.7931
$ws.Range("J135").Value = 166667150
$ws.Range("K135").Value = 6775.1379
$ws.Range("L135").Value = 1500004350
$ws.Range("M135").Value = -4240.1379
$ws.Range("N135").Value = -1500009420
$ws.Range("H138").Value = 2661.4658
$ws.Range("I138").Value = 1946.5294
$ws.Range("J138").Value = 3111.611
$ws.Range("K138").Value = 5839.5882
$ws.Range("L138").Value = 9334.832999999999
$ws.Range("M138").Value = -699.5882000000001
$ws.Range("N138").Value = -19614.833

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 25491.732
$ws.Range("I32").Value = 24750.281
$ws.Range("K32").Value = 24750.281
$ws.Range("M32").Value = -24463.281

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H130").Value = 53489
$ws.Range("J130").Value = 53489
$ws.Range("L130").Value = 53489
$ws.Range("N130").Value = -63529

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H134").Value = 3108.4333
$ws.Range("I134").Value = 1528.4375
$ws.Range("K134").Value = 4585.3125
$ws.Range("M134").Value = -2050.3125

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H2").Value = 9563852
$ws.Range("I2").Value = 56.666668
$ws.Range("J2").Value = 17214888
$ws.Range("K2").Value = 340.000008
$ws.Range("L2").Value = 103289328
$ws.Range("M2").Value = -227.000008
$ws.Range("N2").Value = -103289554
$ws.Range("H4").Value = 680.4
$ws.Range("J4").Value = 3002
$ws.Range("L4").Value = 9006
$ws.Range("N4").Value = -9230
$ws.Range("H9").Value = 250000540
$ws.Range("I9").Value = 333333400
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 1000000200
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = -999999976
$ws.Range("N9").Value = -6448
$ws.Range("H16").Value = 300
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H20").Value = 262.5
$ws.Range("I20").Value = 262.5
$ws.Range("K20").Value = 787.5
$ws.Range("M20").Value = -560.5
$ws.Range("H22").Value = 52997.75
$ws.Range("I22").Value = 55245.5
$ws.Range("J22").Value = 50750
$ws.Range("K22").Value = 165736.5
$ws.Range("L22").Value = 152250
$ws.Range("M22").Value = -165567.5
$ws.Range("N22").Value = -152588
$ws.Range("H23").Value = 610.8182
$ws.Range("I23").Value = 501
$ws.Range("J23").Value = 621.8
$ws.Range("K23").Value = 1503
$ws.Range("L23").Value = 1865.4
$ws.Range("M23").Value = -1268
$ws.Range("N23").Value = -2335.4
$ws.Range("H26").Value = 92.8
$ws.Range("I26").Value = 54.666668
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 164.000004
$ws.Range("L26").Value = 450
$ws.Range("M26").Value = 123.999996
$ws.Range("N26").Value = -1026
$ws.Range("H27").Value = 52997.75
$ws.Range("I27").Value = 55245.5
$ws.Range("J27").Value = 50750
$ws.Range("K27").Value = 165736.5
$ws.Range("L27").Value = 152250
$ws.Range("M27").Value = -165634.5
$ws.Range("N27").Value = -152454
$ws.Range("H32").Value = 950
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 950
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2850
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3416
$ws.Range("H34").Value = 1971.4286
$ws.Range("J34").Value = 2955.5557
$ws.Range("L34").Value = 8866.667099999999
$ws.Range("N34").Value = -9034.667099999999
$ws.Range("H46").Value = 2200
$ws.Range("J46").Value = 2200
$ws.Range("L46").Value = 6600
$ws.Range("N46").Value = -6782
$ws.Range("H51").Value = 9478.308000000001
$ws.Range("I51").Value = 13926.625
$ws.Range("J51").Value = 2361
$ws.Range("K51").Value = 41779.875
$ws.Range("L51").Value = 7083
$ws.Range("M51").Value = -41319.875
$ws.Range("N51").Value = -8003
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2314
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 2098.4
$ws.Range("I64").Value = 1623
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 4869
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -4599
$ws.Range("N64").Value = -12540
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 9000
$ws.Range("M65").Value = -5568
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 2098.4
$ws.Range("I67").Value = 1623
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 4869
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -3933
$ws.Range("N67").Value = -13872
$ws.Range("H86").Value = 1164.5385
$ws.Range("I86").Value = 533
$ws.Range("J86").Value = 1354
$ws.Range("K86").Value = 1599
$ws.Range("L86").Value = 4062
$ws.Range("M86").Value = -413
$ws.Range("N86").Value = -6434
$ws.Range("H89").Value = 1164.5385
$ws.Range("I89").Value = 533
$ws.Range("J89").Value = 1354
$ws.Range("K89").Value = 4797
$ws.Range("L89").Value = 12186
$ws.Range("M89").Value = 1131
$ws.Range("N89").Value = -24042
$ws.Range("H95").Value = 903174.7
$ws.Range("J95").Value = 1354500
$ws.Range("L95").Value = 4063500
$ws.Range("N95").Value = -4067618
$ws.Range("H96").Value = 5500
$ws.Range("J96").Value = 5500
$ws.Range("L96").Value = 16500
$ws.Range("N96").Value = -20618
$ws.Range("H97").Value = 500
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H104").Value = 3000
$ws.Range("I104").Value = 1000
$ws.Range("J104").Value = 4000
$ws.Range("K104").Value = 3000
$ws.Range("L104").Value = 12000
$ws.Range("M104").Value = -379
$ws.Range("N104").Value = -17242
$ws.Range("H105").Value = 336099.66
$ws.Range("J105").Value = 336099.66
$ws.Range("L105").Value = 1008298.98
$ws.Range("N105").Value = -1013540.98
$ws.Range("H108").Value = 1381.1765
$ws.Range("I108").Value = 1331.6666
$ws.Range("J108").Value = 1500
$ws.Range("K108").Value = 3994.9998
$ws.Range("L108").Value = 4500
$ws.Range("M108").Value = -1114.9998
$ws.Range("N108").Value = -10260
$ws.Range("H112").Value = 2503208.5
$ws.Range("I112").Value = 66667332
$ws.Range("J112").Value = 3307.6624
$ws.Range("K112").Value = 200001996
$ws.Range("L112").Value = 9922.9872
$ws.Range("M112").Value = -200000888
$ws.Range("N112").Value = -12138.9872

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H130").Value = 53992
$ws.Range("J130").Value = 53992
$ws.Range("L130").Value = 53992
$ws.Range("N130").Value = -64032
